$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme
try {
  $tcs.Load("Office Theme")
  Write-Output "load ok"
} catch {
  Write-Output "ERROR: $_"
}
